# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the 59fc2548-... row on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-19 12:33:21"
$wsZhCn.Range("H3").Value = "2016-03-19 12:33:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-19 12:33:24"
$wsDeDe.Range("H3").Value = "2016-03-19 12:33:46"
